$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: "Add border to text on title screen" -> split into "Add "
# and "border to text on title screen", with a _GoBack bookmark
# inserted right after "Add ". (Word keeps only one _GoBack bookmark,
# so adding this one automatically removes the old one that sat after
# "Last level".)
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Add border to text on title screen", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $splitPoint = $rng.Start + 4
    $bmRange = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# ------------------------------------------------------------------
# Change 2: strike-through "- sprite animation for attack"
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("- sprite animation for attack", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Font.StrikeThrough = 1
}

# ------------------------------------------------------------------
# Change 3: strike-through "Zombie level needs to be more complex - Hansong"
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Zombie level needs to be more complex - Hansong", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Font.StrikeThrough = 1
}

# ------------------------------------------------------------------
# Change 4: delete the "Add locked door to player's opening cell - Hansong"
# line (including the line break that introduces it), leaving the line
# break before "Turret lever" intact.
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute([char]0x0B + "Add locked door to player" + [char]0x2019 + "s opening cell - Hansong", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rng.Delete()
}
